# The commit removes the trailing blurb that used to follow the last
# "Requisitos" entry ("LOB1036: Geometria Analítica (Requisito fraco)"):
#   - an empty paragraph
#   - an empty, page-break-before paragraph
#   - the "© 2020 . Contact: ..." copyright paragraph
# leaving the empty paragraph + page-break paragraph that originally
# came after the copyright line as the new tail of the document.

$d = $word.ActiveDocument

# Anchor on content (not a hard-coded paragraph index) so the edit is
# resilient to any earlier, unrelated structural differences.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1036*") {
        $anchor = $p
    }
}

if ($anchor -ne $null) {
    $p1 = $anchor.Next()                 # empty paragraph
    $p2 = $p1.Next()                     # empty, pageBreakBefore paragraph
    $p3 = $p2.Next()                     # the "© 2020 ..." paragraph

    $deleteRange = $d.Range($p1.Range.Start, $p3.Range.End)
    $deleteRange.Delete()
}
